# Add two new chapters to the tracking list, keeping the list sorted
# alphabetically by chapter filename:
#   - "load-balancing.md"  -> belongs right before "manage-content.md"
#   - "multitenancy.md"    -> belongs right before "navigation-and-menus.md"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a blank row for "load-balancing.md" above the current
#     "manage-content.md" row (row 18), then clone the plain (non-highlighted)
#     formatting used by the rest of the table onto the B:D cells of the
#     new row.
[void]$ws.Rows("18").Insert()
[void]$ws.Range("B4:D4").Copy()
[void]$ws.Range("B18:D18").PasteSpecial(-4122)

# --- Insert a blank row for "multitenancy.md" above the current
#     "navigation-and-menus.md" row. Because the previous insert shifted
#     everything down by one, that row is now at index 22.
[void]$ws.Rows("22").Insert()
[void]$ws.Range("B4:D4").Copy()
[void]$ws.Range("B22:D22").PasteSpecial(-4122)

# --- Fill in the chapter names. "multitenancy.md" is written first so it
#     becomes the earlier new shared-string entry, matching the source
#     order, followed by "load-balancing.md".
$ws.Range("A22").Value = "multitenancy.md"
$ws.Range("A18").Value = "load-balancing.md"

# Restore the active selection/cursor position.
[void]$ws.Range("G21").Select()
